$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price / volume figures (and the two coin-ranking row
# reshuffles) for the Sat Sep  9 02:41:56 UTC 2023 cryptos list refresh.
# "Price" (column D) cells are always stored as plain text in this sheet
# (note the thousands-dot formatting, e.g. "25.925.29"), so any new value
# that could be mistaken for a number is written with a leading "'" quote
# prefix and the cell style is immediately reset to Normal afterwards -
# this keeps the cell a plain string (no new number format is left behind)
# exactly like the rest of the untouched text cells on the sheet.

$ws.Range('D2').Value = '25.925.29'
$ws.Range('D3').Value = '1.638.49'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = "'215.29"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('E8').Value = '  -0.49%  '
$ws.Range('D9').Value = "'0.0640"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.38%  '
$ws.Range('D10').Value = "'19.64"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.74%  '
$ws.Range('D11').Value = "'0.0795"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').Value = '1.865.14'
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'4.26"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.594.27'
$ws.Range('E14').Value = '  -2.17%  '
$ws.Range('E15').Value = '  -1.17%  '
$ws.Range('D16').Value = '0.0₃0766'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').Value = '25.936.90'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').Value = "'193.12"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('D22').Value = "'9.93"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('E23').Value = '  -0.98%  '
$ws.Range('D24').Value = "'143.99"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = "'1.79"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.52%  '
$ws.Range('B26').Value = 'BinanceUSD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').Value = "'0.128"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.30%  '
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('D29').Value = "'15.55"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('E30').Value = '  -0.55%  '
$ws.Range('D31').Value = "'0.0504"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.20%  '
$ws.Range('E32').Value = '  -1.57%  '
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('E34').Value = '  -3.65%  '
$ws.Range('E35').Value = '  +1.33%  '
$ws.Range('D36').Value = "'0.902"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.43%  '
$ws.Range('D37').Value = '1.139.70'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('D38').Value = "'0.545"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('E40').Value = '  +0.38%  '
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('E42').Value = '  -3.10%  '
$ws.Range('D43').Value = "'99.42"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.76%  '
$ws.Range('D44').Value = "'0.797"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').Value = '1.774.76'
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('E46').Value = '  +2.32%  '
$ws.Range('D47').Value = "'56.61"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('E48').Value = '  +2.87%  '
$ws.Range('D49').Value = "'1.48"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('D50').Value = "'7.65"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.34%  '
$ws.Range('E51').Value = '  -0.91%  '
